$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "42.189.82"
Set-TextValue "E2" "  -3.15%  "
Set-TextValue "D3" "2.219.13"
Set-TextValue "E3" "  -1.96%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.12%  "
Set-TextValue "D5" "107.65"
Set-TextValue "E5" "  -9.70%  "
Set-TextValue "D6" "297.06"
Set-TextValue "E6" "  +12.13%  "
Set-TextValue "D7" "0.625"
Set-TextValue "E7" "  -2.89%  "
Set-TextValue "E8" "  -0.27%  "
Set-TextValue "D9" "0.598"
Set-TextValue "E9" "  -3.41%  "
Set-TextValue "D10" "43.60"
Set-TextValue "E10" "  -8.08%  "
Set-TextValue "D11" "0.0910"
Set-TextValue "E11" "  -3.33%  "
Set-TextValue "D12" "54.44"
Set-TextValue "E12" "  +0.48%  "
Set-TextValue "E13" "  -4.54%  "
Set-TextValue "D14" "0.975"
Set-TextValue "E14" "  +7.66%  "
Set-TextValue "E15" "  -2.64%  "
Set-TextValue "D16" "14.97"
Set-TextValue "E16" "  -2.19%  "
Set-TextValue "D17" "2.550.92"
Set-TextValue "E17" "  -2.15%  "
Set-TextValue "D18" "2.232.85"
Set-TextValue "E18" "  -1.37%  "
Set-TextValue "D19" "42.252.03"
Set-TextValue "E19" "  -2.96%  "
Set-TextValue "D20" "7.40"
Set-TextValue "E20" "  +7.65%  "
Set-TextValue "E21" "  -4.33%  "
Set-TextValue "D22" "72.29"
Set-TextValue "E22" "  +0.30%  "
Set-TextValue "D23" "3.49"
Set-TextValue "E23" "  +21.71%  "
Set-TextValue "E24" "  -3.87%  "
Set-TextValue "D25" "228.24"
Set-TextValue "E25" "  -2.86%  "
Set-TextValue "D26" "9.03"
Set-TextValue "E26" "  -4.69%  "
Set-TextValue "E27" "  -1.69%  "
Set-TextValue "E28" "  -2.66%  "
Set-TextValue "E29" "  -0.71%  "
Set-TextValue "D30" "38.24"
Set-TextValue "E30" "  -8.00%  "
Set-TextValue "E31" "  -4.75%  "
Set-TextValue "D32" "173.64"
Set-TextValue "E32" "  +1.02%  "
Set-TextValue "D33" "20.94"
Set-TextValue "E33" "  -3.19%  "
Set-TextValue "E34" "  -2.13%  "
Set-TextValue "D35" "5.61"
Set-TextValue "E35" "  -1.60%  "
Set-TextValue "D36" "5.10"
Set-TextValue "E36" "  +11.73%  "
Set-TextValue "D37" "4.36"
Set-TextValue "E37" "  +1.76%  "
Set-TextValue "E38" "  -2.98%  "
Set-TextValue "D39" "0.0369"
Set-TextValue "E39" "  -1.79%  "
Set-TextValue "E40" "  -3.52%  "
Set-TextValue "E41" "  -4.33%  "
Set-TextValue "D42" "71.62"
Set-TextValue "E42" "  -3.36%  "
Set-TextValue "E43" "  -1.82%  "
Set-TextValue "E44" "  -0.11%  "
Set-TextValue "D45" "12.56"
Set-TextValue "E45" "  -9.30%  "
Set-TextValue "E46" "  -4.66%  "
Set-TextValue "E47" "  -6.38%  "
Set-TextValue "E48" "  +4.54%  "
Set-TextValue "D49" "103.30"
Set-TextValue "E49" "  +2.30%  "
Set-TextValue "D50" "8.39"
Set-TextValue "E50" "  -1.43%  "
Set-TextValue "D51" "1.63"
Set-TextValue "E51" "  +5.38%  "
